# faturamento_diario.xlsx update
# Adds one more day of billing data for 07/2025 (day 29) into the daily
# revenue table. The new record is inserted right after the existing
# 07/2025 rows (which end at sheet row 29) and before the 06/2025 block,
# so every row from the old row 30 onward shifts down by one.
#
#   Dia=29  total_venda=20532.76  Mes=7  Ano=2025  Periodo=07/2025

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 30 (and everything below it) down by one row, leaving a
# blank row 30 ready to receive the new record.
$ws.Rows.Item(30).Insert()

# Fill in the newly inserted row with the new day's data.
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 20532.76
$ws.Range("C30").Value = 7
$ws.Range("D30").Value = 2025
$ws.Range("E30").Value = "07/2025"
